$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and header date (08-17 -> 08-18)
$ws.Name = "Through 2022-08-18"
$ws.Range("A9").Value = "August (through 08-18)"

# Update August row (row 9): D..I = 2017..2022
$ws.Range("D9").Value = 46
$ws.Range("E9").Value = 29
$ws.Range("F9").Value = 25
$ws.Range("G9").Value = 115
$ws.Range("H9").Value = 97
$ws.Range("I9").Value = 97

# Update Total row (row 10): D..I = 2017..2022
$ws.Range("D10").Value = 511
$ws.Range("E10").Value = 454
$ws.Range("F10").Value = 329
$ws.Range("G10").Value = 736
$ws.Range("H10").Value = 1007
$ws.Range("I10").Value = 1068
